$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yyyy-mm-dd_CSF_Profile")
$ws.Range("J40:J363").ClearContents()
